# The deck's single Design ("Integral") is re-coloured to the stock
# "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# ThemeColorScheme.Item(n).RGB uses the usual VBA RGB() packing
# (R + G*256 + B*65536), so each hex AABBCC below is converted that way:
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation
$theme = $p.Designs.Item(1).SlideMaster.Theme

$theme.ThemeColorScheme.Item(1).RGB = 0
$theme.ThemeColorScheme.Item(2).RGB = 16777215
$theme.ThemeColorScheme.Item(3).RGB = 6968388
$theme.ThemeColorScheme.Item(4).RGB = 15132391
$theme.ThemeColorScheme.Item(5).RGB = 13998939
$theme.ThemeColorScheme.Item(6).RGB = 3243501
$theme.ThemeColorScheme.Item(7).RGB = 10855845
$theme.ThemeColorScheme.Item(8).RGB = 49407
$theme.ThemeColorScheme.Item(9).RGB = 12874308
$theme.ThemeColorScheme.Item(10).RGB = 4697456
$theme.ThemeColorScheme.Item(11).RGB = 12673797
$theme.ThemeColorScheme.Item(12).RGB = 7491477
